$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.626.68'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '2.245.29'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = '306.15'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').Value = '95.05'
$ws.Range('E6').Value = '  -1.96%  '
$ws.Range('D7').Value = '0.572'
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('D10').Value = '35.23'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = '2.586.99'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = '2.237.75'
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = '0.835'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').Value = '44.422.83'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('E19').Value = '  -3.08%  '
$ws.Range('D20').Value = '11.83'
$ws.Range('E20').Value = '  -2.92%  '
$ws.Range('D21').Value = '6.21'
$ws.Range('E21').Value = '  -2.94%  '
$ws.Range('D22').Value = '65.34'
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').Value = '237.37'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  -1.37%  '
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('E27').Value = '  +6.03%  '
$ws.Range('D28').Value = '9.78'
$ws.Range('E28').Value = '  -2.22%  '
$ws.Range('D29').Value = '37.21'
$ws.Range('E29').Value = '  -4.41%  '
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = '149.78'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('D33').Value = '0.0787'
$ws.Range('E33').Value = '  -1.83%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -2.79%  '
$ws.Range('E36').Value = '  +1.38%  '
$ws.Range('E37').Value = '  -1.45%  '
$ws.Range('E38').Value = '  +4.91%  '
$ws.Range('D39').Value = '15.30'
$ws.Range('E39').Value = '  +4.21%  '
$ws.Range('D40').Value = '3.39'
$ws.Range('E40').Value = '  -6.35%  '
$ws.Range('E41').Value = '  -2.30%  '
$ws.Range('D42').Value = '0.0299'
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D44').Value = '1.810.94'
$ws.Range('E44').Value = '  +3.17%  '
$ws.Range('E45').Value = '  +12.13%  '
$ws.Range('D46').Value = '81.60'
$ws.Range('E46').Value = '  -1.98%  '
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').Value = '98.53'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('E49').Value = '  -2.81%  '
$ws.Range('D50').Value = '68.70'
$ws.Range('E50').Value = '  +0.69%  '
$ws.Range('D51').Value = '54.15'
$ws.Range('E51').Value = '  -1.55%  '
